$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Add the three new character styles used by the edited paragraphs.
#    wdStyleTypeCharacter = 2
# ---------------------------------------------------------------------
$ganStyle = $d.Styles.Add("GaNStyle", 2)
$ganStyle.Font.Name = "Calibri"
$ganStyle.Font.Size = 14

$ganParagraph = $d.Styles.Add("GaNParagraph", 2)
$ganParagraph.Font.Name = "Calibri"
$ganParagraph.Font.Size = 10

$ganLinks = $d.Styles.Add("GaNLinks", 2)
$ganLinks.Font.Name = "Calibri"
$ganLinks.Font.Bold = $true
$ganLinks.Font.Color = 8388608
$ganLinks.Font.Size = 9.5
$ganLinks.Font.Underline = 1

# ---------------------------------------------------------------------
# 2. Apply GaNStyle to every "Waktu Kampanye 2022 ..." run (4 places).
# ---------------------------------------------------------------------
$rng = $d.Content
while ($rng.Find.Execute("Waktu Kampanye 2022 untuk konstelasi Perseus: 16-25 Januari, 7-16 November, 6-15 Desember", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
}

# ---------------------------------------------------------------------
# 3. Apply GaNParagraph to the "Anda sedang berpartisipasi ..." run.
# ---------------------------------------------------------------------
$rng = $d.Content
if ($rng.Find.Execute("Anda sedang berpartisipasi dalam kampanye global pengamatan dan pencatatan penampakan bintang paling redup untuk pengukuran tingkat polusi cahaya di suatu lokasi. Melalui pengamatan dan identifikasi  konstelasi Perseus di langit malam dan membandingkannya dengan peta bintang, masyarakat di seluruh dunia dapat mengetahui dan mempelajari seberapa besar kontribusi cahaya di lingkungannya terhadap polusi cahaya. Kontribusi data anda pada basis data online akan membantu mendokumentasikan langit malam yang tampak di berbagai lokasi.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNParagraph"
}

# ---------------------------------------------------------------------
# 4. Apply GaNLinks to the "Peta di dokumen ini disiapkan oleh ..." run.
# ---------------------------------------------------------------------
$rng = $d.Content
if ($rng.Find.Execute("Peta di dokumen ini disiapkan oleh Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNLinks"
}

Write-Host "Styles created and applied."
